$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (border/font/alignment) of the last existing A-column data cell (A205)
# onto each new A-column cell so the appended rows keep the same "s=1" formatting as the rest of the column.
$srcStyle = $ws.Range("A205")

$srcStyle.Copy()
$ws.Range("A206").PasteSpecial(-4122)
$ws.Range("A206").Value = 204
$ws.Range("B206").Value = 0.5242630385487528

$srcStyle.Copy()
$ws.Range("A207").PasteSpecial(-4122)
$ws.Range("A207").Value = 205
$ws.Range("B207").Value = 0.5255823541537827

$srcStyle.Copy()
$ws.Range("A208").PasteSpecial(-4122)
$ws.Range("A208").Value = 206
$ws.Range("B208").Value = 0.3921995464852608

$srcStyle.Copy()
$ws.Range("A209").PasteSpecial(-4122)
$ws.Range("A209").Value = 207
$ws.Range("B209").Value = 0.490400604686319

$srcStyle.Copy()
$ws.Range("A210").PasteSpecial(-4122)
$ws.Range("A210").Value = 208
$ws.Range("B210").Value = 0.3573696145124716

$srcStyle.Copy()
$ws.Range("A211").PasteSpecial(-4122)
$ws.Range("A211").Value = 209
$ws.Range("B211").Value = 0.5781665046971169

$srcStyle.Copy()
$ws.Range("A212").PasteSpecial(-4122)
$ws.Range("A212").Value = 210
$ws.Range("B212").Value = 0.4698412698412698

$srcStyle.Copy()
$ws.Range("A213").PasteSpecial(-4122)
$ws.Range("A213").Value = 211
$ws.Range("B213").Value = 0.6534240362811791

$srcStyle.Copy()
$ws.Range("A214").PasteSpecial(-4122)
$ws.Range("A214").Value = 212
$ws.Range("B214").Value = 0.6839002267573697

$srcStyle.Copy()
$ws.Range("A215").PasteSpecial(-4122)
$ws.Range("A215").Value = 213
$ws.Range("B215").Value = 0.509750566893424

$srcStyle.Copy()
$ws.Range("A216").PasteSpecial(-4122)
$ws.Range("A216").Value = 214
$ws.Range("B216").Value = 0.509750566893424

$srcStyle.Copy()
$ws.Range("A217").PasteSpecial(-4122)
$ws.Range("A217").Value = 215
$ws.Range("B217").Value = 0.509750566893424

$excel.CutCopyMode = 0

Write-Host "Appended rows 206:217"